{"js": "// Update the 25 division-problem answers that live in the first table of\n// the document. The table has a repeating 4-row block pattern: one row of\n// 5 filled answer cells followed by 3 empty spacer rows. We walk every\n// cell in document order and, whenever its text matches the next expected\n// \"before\" value, replace just that text run's content (via a Range over\n// the paragraph) so paragraph/run formatting (font, size, alignment) is\n// left completely untouched \u2014 only the visible digits change, exactly as\n// in the source diff.\n\nconst replacements = [\n  [\"92\u00f76=15, 2\", \"46\u00f73=15, 1\"],\n  [\"80\u00f76=13, 2\", \"91\u00f73=30, 1\"],\n  [\"37\u00f72=18, 1\", \"39\u00f75=7, 4\"],\n  [\"59\u00f73=19, 2\", \"83\u00f79=9, 2\"],\n  [\"43\u00f78=5, 3\", \"48\u00f75=9, 3\"],\n  [\"97\u00f79=10, 7\", \"53\u00f73=17, 2\"],\n  [\"15\u00f76=2, 3\", \"90\u00f74=22, 2\"],\n  [\"20\u00f75=4, 0\", \"84\u00f74=21, 0\"],\n  [\"66\u00f75=13, 1\", \"29\u00f79=3, 2\"],\n  [\"14\u00f73=4, 2\", \"79\u00f78=9, 7\"],\n  [\"47\u00f76=7, 5\", \"81\u00f72=40, 1\"],\n  [\"36\u00f72=18, 0\", \"13\u00f74=3, 1\"],\n  [\"67\u00f76=11, 1\", \"59\u00f78=7, 3\"],\n  [\"99\u00f75=19, 4\", \"89\u00f78=11, 1\"],\n  [\"65\u00f77=9, 2\", \"58\u00f78=7, 2\"],\n  [\"52\u00f77=7, 3\", \"92\u00f74=23, 0\"],\n  [\"79\u00f75=15, 4\", \"24\u00f75=4, 4\"],\n  [\"43\u00f78=5, 3\", \"69\u00f76=11, 3\"],\n  [\"53\u00f74=13, 1\", \"84\u00f76=14, 0\"],\n  [\"36\u00f78=4, 4\", \"68\u00f73=22, 2\"],\n  [\"97\u00f78=12, 1\", \"98\u00f72=49, 0\"],\n  [\"40\u00f73=13, 1\", \"36\u00f75=7, 1\"],\n  [\"49\u00f77=7, 0\", \"63\u00f72=31, 1\"],\n  [\"64\u00f75=12, 4\", \"56\u00f79=6, 2\"],\n  [\"95\u00f74=23, 3\", \"22\u00f79=2, 4\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nfor (const row of rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\nlet idx = 0;\nfor (const row of rows.items) {\n  if (idx >= replacements.length) break;\n  for (const cell of row.cells.items) {\n    if (idx >= replacements.length) break;\n\n    const para = cell.body.paragraphs.getFirst();\n    const range = para.getRange();\n    range.load(\"text\");\n    await context.sync();\n\n    const [oldText, newText] = replacements[idx];\n    if (range.text === oldText) {\n      range.insertText(newText, Word.InsertLocation.replace);\n      idx++;\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the 25 division-problem answers that live in the first table of\n# the document. The table has a repeating 4-row block pattern: one row of\n# 5 filled answer cells followed by 3 empty spacer rows. We walk every\n# cell of the table in row-major order and, whenever its (trimmed) text\n# matches the next expected \"before\" value, assign the new text straight\n# onto the cell's Range. Word keeps the existing run/paragraph formatting\n# (font, size, alignment) intact when Range.Text is reassigned like this,\n# so only the visible digits change, exactly as in the source diff.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$olds = @(\n  \"92\u00f76=15, 2\", \"80\u00f76=13, 2\", \"37\u00f72=18, 1\", \"59\u00f73=19, 2\", \"43\u00f78=5, 3\",\n  \"97\u00f79=10, 7\", \"15\u00f76=2, 3\", \"20\u00f75=4, 0\", \"66\u00f75=13, 1\", \"14\u00f73=4, 2\",\n  \"47\u00f76=7, 5\", \"36\u00f72=18, 0\", \"67\u00f76=11, 1\", \"99\u00f75=19, 4\", \"65\u00f77=9, 2\",\n  \"52\u00f77=7, 3\", \"79\u00f75=15, 4\", \"43\u00f78=5, 3\", \"53\u00f74=13, 1\", \"36\u00f78=4, 4\",\n  \"97\u00f78=12, 1\", \"40\u00f73=13, 1\", \"49\u00f77=7, 0\", \"64\u00f75=12, 4\", \"95\u00f74=23, 3\"\n)\n$news = @(\n  \"46\u00f73=15, 1\", \"91\u00f73=30, 1\", \"39\u00f75=7, 4\", \"83\u00f79=9, 2\", \"48\u00f75=9, 3\",\n  \"53\u00f73=17, 2\", \"90\u00f74=22, 2\", \"84\u00f74=21, 0\", \"29\u00f79=3, 2\", \"79\u00f78=9, 7\",\n  \"81\u00f72=40, 1\", \"13\u00f74=3, 1\", \"59\u00f78=7, 3\", \"89\u00f78=11, 1\", \"58\u00f78=7, 2\",\n  \"92\u00f74=23, 0\", \"24\u00f75=4, 4\", \"69\u00f76=11, 3\", \"84\u00f76=14, 0\", \"68\u00f73=22, 2\",\n  \"98\u00f72=49, 0\", \"36\u00f75=7, 1\", \"63\u00f72=31, 1\", \"56\u00f79=6, 2\", \"22\u00f79=2, 4\"\n)\n\n$idx = 0\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n  for ($c = 1; $c -le $t.Columns.Count; $c++) {\n    if ($idx -ge $olds.Length) { break }\n\n    $cell = $t.Cell($r, $c)\n    $txt = $cell.Range.Text.TrimEnd([char]13, [char]7)\n\n    if ($txt -eq $olds[$idx]) {\n      $cell.Range.Text = $news[$idx]\n      $idx++\n    }\n  }\n}\n"}
